# Applies the cryptos.xlsx price/volume/coin update described by the commit
# "Updated cryptos list on Sun Sep 24 19:58:28 UTC 2023 with GitHub Actions".
#
# All touched cells in this sheet are plain-text cells (t="inlineStr" in the
# original OOXML) even though many of their contents look like numbers
# (e.g. "211.00", "7.13", "0.0511"). Assigning such a string straight to
# Range.Value lets Excel's normal type-inference turn it into a real number
# (and drop formatting like trailing zeros), which would not match the
# source data. To preserve the literal text we briefly mark the cell as
# Text ("@") before assigning the value, then clear the format again so we
# do not leave a stray number-format style behind on a cell that originally
# had the default style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($Sheet, $Address, $Text)
    $cell = $Sheet.Range($Address)
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.ClearFormats()
}

# Row 2
Set-TextCell $ws "D2" "26.651.88"
Set-TextCell $ws "E2" "  -0.17%  "

# Row 3
Set-TextCell $ws "D3" "1.597.87"
Set-TextCell $ws "E3" "  -0.13%  "

# Row 4
Set-TextCell $ws "E4" "  +0.37%  "

# Row 5
Set-TextCell $ws "D5" "211.00"
Set-TextCell $ws "E5" "  -0.19%  "

# Row 6
Set-TextCell $ws "E6" "  -0.38%  "

# Row 7
Set-TextCell $ws "E7" "  +0.38%  "

# Row 8
Set-TextCell $ws "E8" "  -0.46%  "

# Row 9
Set-TextCell $ws "E9" "  -0.52%  "

# Row 10
Set-TextCell $ws "E10" "  +0.24%  "

# Row 11
Set-TextCell $ws "E11" "  +0.48%  "

# Row 12
Set-TextCell $ws "D12" "1.821.49"
Set-TextCell $ws "E12" "  -0.18%  "

# Row 13
Set-TextCell $ws "D13" "1.597.12"
Set-TextCell $ws "E13" "  +0.42%  "

# Row 14
Set-TextCell $ws "E14" "  +0.03%  "

# Row 16
Set-TextCell $ws "E16" "  -0.92%  "

# Row 17
Set-TextCell $ws "D17" "26.630.01"

# Row 19
Set-TextCell $ws "E19" "  +0.38%  "

# Row 20
Set-TextCell $ws "D20" "208.55"
Set-TextCell $ws "E20" "  -0.80%  "

# Row 21
Set-TextCell $ws "D21" "7.13"
Set-TextCell $ws "E21" "  -1.40%  "

# Row 22
Set-TextCell $ws "E22" "  -0.09%  "

# Row 23
Set-TextCell $ws "E23" "  -2.73%  "

# Row 24
Set-TextCell $ws "D24" "8.96"
Set-TextCell $ws "E24" "  +0.22%  "

# Row 25
Set-TextCell $ws "D25" "143.90"
Set-TextCell $ws "E25" "  +0.55%  "

# Row 26
Set-TextCell $ws "E26" "  +0.66%  "

# Row 27
Set-TextCell $ws "E27" "  +0.23%  "

# Row 28
Set-TextCell $ws "E28" "  -0.73%  "

# Row 30
Set-TextCell $ws "E30" "  -2.29%  "

# Row 31
Set-TextCell $ws "E31" "  -0.26%  "

# Row 32
Set-TextCell $ws "E32" "  -0.37%  "

# Row 33
Set-TextCell $ws "E33" "  -0.02%  "

# Row 34
Set-TextCell $ws "E34" "  +19.71%  "

# Row 35
Set-TextCell $ws "D35" "1.278.00"
Set-TextCell $ws "E35" "  -0.93%  "

# Row 36
Set-TextCell $ws "D36" "2.49"
Set-TextCell $ws "E36" "  +1.08%  "

# Row 37
Set-TextCell $ws "B37" "LidoDAOToken"
Set-TextCell $ws "C37" "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextCell $ws "D37" "1.49"
Set-TextCell $ws "E37" "  -0.78%  "

# Row 38
Set-TextCell $ws "B38" "ImmutableX"
Set-TextCell $ws "C38" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextCell $ws "D38" "0.598"
Set-TextCell $ws "E38" "  -3.65%  "

# Row 39
Set-TextCell $ws "E39" "  -2.04%  "

# Row 40
Set-TextCell $ws "E40" "  -0.43%  "

# Row 41
Set-TextCell $ws "B41" "FraxShare"
Set-TextCell $ws "C41" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextCell $ws "D41" "5.42"
Set-TextCell $ws "E41" "  -0.19%  "

# Row 42
Set-TextCell $ws "B42" "TrustWalletToken"
Set-TextCell $ws "C42" "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextCell $ws "D42" "0.777"
Set-TextCell $ws "E42" "  -1.06%  "

# Row 43
Set-TextCell $ws "B43" "MXToken"
Set-TextCell $ws "C43" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextCell $ws "D43" "2.16"
Set-TextCell $ws "E43" "  -0.93%  "

# Row 45
Set-TextCell $ws "D45" "1.732.93"
Set-TextCell $ws "E45" "  -0.06%  "

# Row 46
Set-TextCell $ws "D46" "90.20"
Set-TextCell $ws "E46" "  -1.14%  "

# Row 47
Set-TextCell $ws "D47" "1.56"
Set-TextCell $ws "E47" "  -0.70%  "

# Row 48
Set-TextCell $ws "E48" "  -2.32%  "

# Row 49
Set-TextCell $ws "E49" "  +2.31%  "

# Row 50
Set-TextCell $ws "D50" "0.0511"
Set-TextCell $ws "E50" "  +0.49%  "

# Row 51
Set-TextCell $ws "B51" "USDD"
Set-TextCell $ws "C51" "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
Set-TextCell $ws "D51" "1.01"
Set-TextCell $ws "E51" "  +0.29%  "
